# "add mouser link to fuse"
# Adds a new part-list row (F1-3 / Fuse) with a Mouser hyperlink in column C,
# matching rows 2-11's existing pattern (col A = designator, col B = name,
# col C = hyperlinked Mouser/Farnell link).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fuseLink = "https://www.mouser.de/ProductDetail/Bel-Fuse/C1T5?qs=GtFly9OVs8%2FF1GxRAaUoTA%3D%3D"

# Column C: hyperlink cell, display text = the URL itself (same as the other
# link cells in this sheet), styled with the built-in Hyperlink cell style.
$ws.Hyperlinks.Add($ws.Range("C12"), $fuseLink, "", "", $fuseLink)
$ws.Range("C12").Style = "Hyperlink"

# Column B: part name, Column A: schematic designator.
$ws.Range("B12").Value = "Fuse"
$ws.Range("A12").Value = "F1-3"

# Leave the selection where the author ended up after entering the row.
[void]$ws.Range("B16").Select()
